# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45204 (2023-10-05) to 45207 (2023-10-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C56").Value = 45207
